$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("L13").Value = 24666.334
$ws.Range("N13").Value = -25004.334
$ws.Range("H13").Value = 24749.75
$ws.Range("J13").Value = 24666.334
$ws.Range("L33").Value = 4991.3335
$ws.Range("M33").Value = -77.55883999999998
$ws.Range("N33").Value = -5449.3335
$ws.Range("K33").Value = 306.55884
$ws.Range("H33").Value = 686.4054
$ws.Range("I33").Value = 306.55884
$ws.Range("J33").Value = 4991.3335
$ws.Range("L41").Value = 454.5
$ws.Range("J41").Value = 454.5
$ws.Range("H41").Value = 838.38464
$ws.Range("I41").Value = 1009
$ws.Range("K41").Value = 1009
$ws.Range("N41").Value = -1334.5
$ws.Range("M41").Value = -569
$ws.Range("I53").Value = 181.25
$ws.Range("H53").Value = 165
$ws.Range("M53").Value = 455.75
$ws.Range("K53").Value = 181.25
$ws.Range("M61").Value = -2976.5
$ws.Range("K61").Value = 3148.5
$ws.Range("I61").Value = 1049.5
$ws.Range("H61").Value = 1049.5
$ws.Range("J87").Value = 131563.5
$ws.Range("H87").Value = 131563.5
$ws.Range("L87").Value = 131563.5
$ws.Range("N87").Value = -134059.5
$ws.Range("H90").Value = 131563.5
$ws.Range("N90").Value = -407170.5
$ws.Range("J90").Value = 131563.5
$ws.Range("L90").Value = 394690.5
$ws.Range("I100").Value = 1909
$ws.Range("L100").Value = 2310.9443
$ws.Range("J100").Value = 2310.9443
$ws.Range("K100").Value = 1909
$ws.Range("H100").Value = 2150.1667
$ws.Range("N100").Value = -3392.9443
$ws.Range("M100").Value = -1368
$ws.Range("M112").Value = -3149.75
$ws.Range("H112").Value = 5456.3076
$ws.Range("K112").Value = 4257.75
$ws.Range("N112").Value = -20786.9552
$ws.Range("J112").Value = 6190.3184
$ws.Range("L112").Value = 18570.9552
$ws.Range("I112").Value = 1419.25
$ws.Range("N113").Value = -11973
$ws.Range("L113").Value = 5465
$ws.Range("H113").Value = 6882.3335
$ws.Range("J113").Value = 5465
$ws.Range("H129").Value = 1404.3334
$ws.Range("J129").Value = 4000
$ws.Range("N129").Value = -22000
$ws.Range("L129").Value = 12000
$ws.Range("N130").Value = ""
$ws.Range("L130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("H130").Value = 0
$ws.Range("K138").Value = 3798.15
$ws.Range("N138").Value = -22629.179
$ws.Range("L138").Value = 12349.179
$ws.Range("M138").Value = 1341.85
$ws.Range("H138").Value = 2928.75
$ws.Range("I138").Value = 1266.05
$ws.Range("J138").Value = 4116.393

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J2").Value = 4999.5
$ws.Range("H2").Value = 2494.6
$ws.Range("N2").Value = -5225.5
$ws.Range("L2").Value = 4999.5
$ws.Range("N5").Value = -5222.5
$ws.Range("H5").Value = 4335.273
$ws.Range("L5").Value = 4998.5
$ws.Range("J5").Value = 4998.5
$ws.Range("I5").Value = 4187.8887
$ws.Range("K5").Value = 4187.8887
$ws.Range("M5").Value = -4075.8887
$ws.Range("L10").Value = 1000
$ws.Range("N10").Value = -1340
$ws.Range("H10").Value = 1000
$ws.Range("J10").Value = 1000
$ws.Range("L32").Value = 31499.25
$ws.Range("H32").Value = 7365.7856
$ws.Range("N32").Value = -32073.25
$ws.Range("J32").Value = 31499.25
$ws.Range("K45").Value = 35638
$ws.Range("M45").Value = -35261
$ws.Range("N45").Value = -2424.5
$ws.Range("L45").Value = 1670.5
$ws.Range("H45").Value = 27146.125
$ws.Range("I45").Value = 35638
$ws.Range("J45").Value = 1670.5
$ws.Range("H46").Value = 9015.4
$ws.Range("N46").Value = -10414.5
$ws.Range("J46").Value = 9776.5
$ws.Range("L46").Value = 9776.5
$ws.Range("H86").Value = 10000
$ws.Range("I86").Value = 10000
$ws.Range("K86").Value = 10000
$ws.Range("M86").Value = -8814
$ws.Range("I89").Value = 10000
$ws.Range("K89").Value = 30000
$ws.Range("H89").Value = 10000
$ws.Range("M89").Value = -24072
$ws.Range("N116").Value = -9587.5
$ws.Range("J116").Value = 4999.5
$ws.Range("L116").Value = 4999.5
$ws.Range("H116").Value = 2494.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J3").Value = 4999.5
$ws.Range("H3").Value = 2494.6
$ws.Range("N3").Value = -5227.5
$ws.Range("L3").Value = 4999.5
$ws.Range("I4").Value = 4187.8887
$ws.Range("K4").Value = 4187.8887
$ws.Range("N4").Value = -5228.5
$ws.Range("H4").Value = 4335.273
$ws.Range("L4").Value = 4998.5
$ws.Range("J4").Value = 4998.5
$ws.Range("M4").Value = -4072.8887
$ws.Range("L80").Value = 2001.25
$ws.Range("H80").Value = 2057.7778
$ws.Range("J80").Value = 2001.25
$ws.Range("N80").Value = -3997.25
$ws.Range("J83").Value = 2001.25
$ws.Range("H83").Value = 2057.7778
$ws.Range("N83").Value = -19990.25
$ws.Range("L83").Value = 10006.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L10").Value = 10008
$ws.Range("M10").Value = -689.25
$ws.Range("N10").Value = -10286
$ws.Range("H10").Value = 2664.2
$ws.Range("I10").Value = 828.25
$ws.Range("J10").Value = 10008
$ws.Range("K10").Value = 828.25
$ws.Range("H16").Value = 38148.625
$ws.Range("I16").Value = 25862.834
$ws.Range("K16").Value = 25862.834
$ws.Range("M16").Value = -25575.834
$ws.Range("L31").Value = 3778.3333
$ws.Range("J31").Value = 3778.3333
$ws.Range("H31").Value = 6927.1333
$ws.Range("N31").Value = -4368.3333
$ws.Range("J34").Value = 3778.3333
$ws.Range("N34").Value = -4182.3333
$ws.Range("H34").Value = 6927.1333
$ws.Range("L34").Value = 3778.3333
$ws.Range("K58").Value = 2098
$ws.Range("I58").Value = 2098
$ws.Range("H58").Value = 2098
$ws.Range("M58").Value = -1895
$ws.Range("J58").Value = 0
$ws.Range("N58").Value = ""
$ws.Range("L58").Value = 0
$ws.Range("K107").Value = 2610.6155
$ws.Range("H107").Value = 2595.5334
$ws.Range("I107").Value = 2610.6155
$ws.Range("M107").Value = -690.6154999999999
$ws.Range("M113").Value = -23692.834
$ws.Range("H113").Value = 38148.625
$ws.Range("K113").Value = 25862.834
$ws.Range("I113").Value = 25862.834
$ws.Range("K132").Value = 15685.9995
$ws.Range("I132").Value = 5228.6665
$ws.Range("M132").Value = -13155.9995
$ws.Range("H132").Value = 5089.375
$ws.Range("H136").Value = 2098
$ws.Range("N136").Value = ""
$ws.Range("L136").Value = 0
$ws.Range("K136").Value = 6294
$ws.Range("M136").Value = -3744
$ws.Range("J136").Value = 0
$ws.Range("I136").Value = 2098

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J2").Value = 356.66666
$ws.Range("H2").Value = 1013.6316
$ws.Range("N2").Value = -2365.99996
$ws.Range("L2").Value = 2139.99996
$ws.Range("K6").Value = 204
$ws.Range("M6").Value = -91
$ws.Range("I6").Value = 68
$ws.Range("L6").Value = 0
$ws.Range("H6").Value = 68
$ws.Range("J6").Value = 0
$ws.Range("N6").Value = ""
$ws.Range("J34").Value = 1194
$ws.Range("N34").Value = -3750
$ws.Range("H34").Value = 1341.1666
$ws.Range("L34").Value = 3582
$ws.Range("H39").Value = 5268
$ws.Range("J39").Value = 6070.1113
$ws.Range("N39").Value = -18798.3339
$ws.Range("L39").Value = 18210.3339
$ws.Range("J54").Value = 7499.5
$ws.Range("N54").Value = -23616.5
$ws.Range("L54").Value = 22498.5
$ws.Range("H54").Value = 7499.5
$ws.Range("H55").Value = 5097.125
$ws.Range("L55").Value = 15715.0005
$ws.Range("N55").Value = -16069.0005
$ws.Range("J55").Value = 5238.3335
$ws.Range("L59").Value = 19498.5
$ws.Range("N59").Value = -20578.5
$ws.Range("J59").Value = 6499.5
$ws.Range("M59").Value = -3963.75
$ws.Range("H59").Value = 3167.3333
$ws.Range("K59").Value = 4503.75
$ws.Range("I59").Value = 1501.25
$ws.Range("L64").Value = 35500.875
$ws.Range("H64").Value = 10627.923
$ws.Range("K64").Value = 26096.4
$ws.Range("I64").Value = 8698.799999999999
$ws.Range("M64").Value = -25826.4
$ws.Range("N64").Value = -36040.875
$ws.Range("J64").Value = 11833.625
$ws.Range("H67").Value = 10627.923
$ws.Range("N67").Value = -37372.875
$ws.Range("J67").Value = 11833.625
$ws.Range("M67").Value = -25160.4
$ws.Range("L67").Value = 35500.875
$ws.Range("I67").Value = 8698.799999999999
$ws.Range("K67").Value = 26096.4
$ws.Range("N131").Value = -17887.8297
$ws.Range("L131").Value = 7807.8297
$ws.Range("H131").Value = 2263.4814
$ws.Range("J131").Value = 2602.6099
$ws.Range("N137").Value = -26913
$ws.Range("L137").Value = 16713
$ws.Range("J137").Value = 5571
$ws.Range("H137").Value = 5586.154

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J94").Value = 41000
$ws.Range("N94").Value = -42352
$ws.Range("L94").Value = 41000
$ws.Range("H94").Value = 41000
$ws.Range("K132").Value = 18549.4995
$ws.Range("I132").Value = 6183.1665
$ws.Range("N132").Value = -15245
$ws.Range("J132").Value = 3395
$ws.Range("M132").Value = -16019.4995
$ws.Range("L132").Value = 10185
$ws.Range("H132").Value = 5784.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M22").Value = -2153.3572
$ws.Range("L22").Value = 2908.182
$ws.Range("I22").Value = 2448.3572
$ws.Range("K22").Value = 2448.3572
$ws.Range("H22").Value = 2650.68
$ws.Range("J22").Value = 2908.182
$ws.Range("N22").Value = -3498.182
$ws.Range("J27").Value = 2908.182
$ws.Range("L27").Value = 2908.182
$ws.Range("N27").Value = -3122.182
$ws.Range("M27").Value = -2341.3572
$ws.Range("H27").Value = 2650.68
$ws.Range("I27").Value = 2448.3572
$ws.Range("K27").Value = 2448.3572
$ws.Range("J40").Value = 22595
$ws.Range("I40").Value = 30373.75
$ws.Range("L40").Value = 22595
$ws.Range("H40").Value = 26052.223
$ws.Range("N40").Value = -22867
$ws.Range("K40").Value = 30373.75
$ws.Range("M40").Value = -30237.75
$ws.Range("I46").Value = 5749.25
$ws.Range("H46").Value = 7624.625
$ws.Range("M46").Value = -5561.25
$ws.Range("K46").Value = 5749.25
$ws.Range("H55").Value = 299.27274
$ws.Range("L55").Value = 148.83333
$ws.Range("I55").Value = 479.8
$ws.Range("K55").Value = 479.8
$ws.Range("N55").Value = -494.83333
$ws.Range("M55").Value = -306.8
$ws.Range("J55").Value = 148.83333
$ws.Range("M61").Value = -1488
$ws.Range("K61").Value = 1690
$ws.Range("I61").Value = 1690
$ws.Range("H61").Value = 1690
$ws.Range("I82").Value = 2170.1333
$ws.Range("M82").Value = -1809.1333
$ws.Range("H82").Value = 1758.5
$ws.Range("K82").Value = 2170.1333
$ws.Range("M85").Value = -922.1333
$ws.Range("H85").Value = 1758.5
$ws.Range("I85").Value = 2170.1333
$ws.Range("K85").Value = 2170.1333
$ws.Range("I93").Value = 1102.2941
$ws.Range("N93").Value = -5564
$ws.Range("J93").Value = 3068
$ws.Range("M93").Value = 145.7058999999999
$ws.Range("K93").Value = 1102.2941
$ws.Range("H93").Value = 1397.15
$ws.Range("L93").Value = 3068
$ws.Range("M113").Value = 480
$ws.Range("H113").Value = 1690
$ws.Range("K113").Value = 1690
$ws.Range("I113").Value = 1690
$ws.Range("K122").Value = 54177
$ws.Range("L122").Value = 0
$ws.Range("I122").Value = 18059
$ws.Range("J122").Value = 0
$ws.Range("H122").Value = 18059
$ws.Range("N122").Value = ""
$ws.Range("M122").Value = -51727
$ws.Range("K132").Value = 15525
$ws.Range("I132").Value = 5175
$ws.Range("M132").Value = -12995
$ws.Range("J132").Value = 3999.875
$ws.Range("N132").Value = -17059.625
$ws.Range("L132").Value = 11999.625
$ws.Range("H132").Value = 4939.975
$ws.Range("H136").Value = 1270.2858
$ws.Range("K136").Value = 3896.5002
$ws.Range("M136").Value = -1346.5002
$ws.Range("I136").Value = 1298.8334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K122").Value = 8120.625
$ws.Range("L122").Value = 7660.5
$ws.Range("I122").Value = 2706.875
$ws.Range("J122").Value = 2553.5
$ws.Range("H122").Value = 2655.75
$ws.Range("N122").Value = -12560.5
$ws.Range("M122").Value = -5670.625
$ws.Range("H136").Value = 2528.6538
$ws.Range("K136").Value = 5227.950000000001
$ws.Range("M136").Value = -2677.950000000001
$ws.Range("I136").Value = 1742.65
